$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append after the last existing row (301): row, date-serial, B, C, D
$newRows = @(
  @(302, 44376, 0, 0, 0),
  @(303, 44377, 0, 0, 0),
  @(304, 44378, 0, 0, 0),
  @(305, 44379, 0, 0, 0),
  @(306, 44380, 0, 0, 0),
  @(307, 44381, 0, 0, 0),
  @(308, 44382, 0, 0, 0),
  @(309, 44383, 0, 0, 0),
  @(310, 44384, 0, 0, 0),
  @(311, 44385, 0, 0, 0),
  @(312, 44386, 0, 0, 0),
  @(313, 44387, 0, 0, 0),
  @(314, 44388, 0, 0, 0),
  @(315, 44389, 0, 0, 0),
  @(316, 44390, 0, 0, 0),
  @(317, 44391, 0, 0, 0),
  @(318, 44392, 0, 0, 0),
  @(319, 44393, 0, 0, 0),
  @(320, 44394, 0, 0, 0),
  @(321, 44395, 0, 0, 0),
  @(322, 44396, 1, 1, 28.87669650591972),
  @(323, 44397, 0, 1, 28.87669650591972),
  @(324, 44398, 0, 1, 28.87669650591972),
  @(325, 44399, 0, 1, 28.87669650591972),
  @(326, 44400, 0, 1, 28.87669650591972),
  @(327, 44401, 0, 1, 28.87669650591972),
  @(328, 44402, 0, 1, 28.87669650591972)
)

$lastRow = 301
$firstNew = 302
$lastNew = 328

# Carry the formatting of the last existing row (date style on column A, plain
# numbers on B:D) down onto the newly appended rows before filling in values.
$ws.Range("A$lastRow`:D$lastRow").Copy()
$ws.Range("A$firstNew`:D$lastNew").PasteSpecial(-4122)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
}
